$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NroSiniestro value for the first "Mediación" row (F2):
# old claim number 0420172010219 -> new claim number 0420172010448
# (leading apostrophe keeps it text, preserving the leading zero / quote-prefix style)
$ws.Range("F2").Value = "'0420172010448"

# Move the active cell selection from F7 to H5
$ws.Range("H5").Select()
